# feat: add 2022-Q4 data
#
# - Insert a new worksheet "2022-Q4" right after "总计", populated with the
#   new quarter's fund-holdings data (previously "2022-Q3" occupied this
#   slot; all the other quarter sheets shift right by one position, but
#   their contents are untouched).
# - Prepend a new "2022-Q4" row to the "总计" (totals) summary sheet and
#   push the existing rows down.

$wb = $excel.ActiveWorkbook
$totals = $wb.Worksheets.Item(1)          # "总计"

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $totals)
$q4.Name = "2022-Q4"

# NOTE: worksheet references resolve by *index*, which the Add() above
# just shifted - so the "2022-Q3" sheet (our style donor, used below)
# must be looked up *after* insertion, when it now sits at position 3.
$oldQ3 = $wb.Worksheets.Item(3)           # "2022-Q3" (pre-existing)

# Header row (column B..H)
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Row 2 - fund 010543 (中加科鑫混合A)
# Numeric-looking identifiers/figures ("010543", "1.17", ...) are stored
# as text in the source data, so a leading apostrophe forces text entry
# instead of silently becoming a number (which would also eat the
# leading zero in the fund code).
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'010543"
$q4.Range("B2").Style = "Normal"
$q4.Range("C2").Value = "中加科鑫混合A"
$q4.Range("D2").Value = "'1.17"
$q4.Range("D2").Style = "Normal"
$q4.Range("E2").Value = "'25.19"
$q4.Range("E2").Style = "Normal"
$q4.Range("F2").Value = "'1.17"
$q4.Range("F2").Style = "Normal"
$q4.Range("G2").Value = "'0.0137"
$q4.Range("G2").Style = "Normal"
$q4.Range("H2").Value = 5

# Row 3 - fund 010544 (中加科鑫混合C)
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'010544"
$q4.Range("B3").Style = "Normal"
$q4.Range("C3").Value = "中加科鑫混合C"
$q4.Range("D3").Value = "'0.01"
$q4.Range("D3").Style = "Normal"
$q4.Range("E3").Value = "'25.19"
$q4.Range("E3").Style = "Normal"
$q4.Range("F3").Value = "'1.17"
$q4.Range("F3").Style = "Normal"
$q4.Range("G3").Value = "'0.0001"
$q4.Range("G3").Style = "Normal"
$q4.Range("H3").Value = 5

# Match the header / index-column formatting used by the neighbouring
# quarter sheets (bold, centred, bordered "s=2" style).
$oldQ3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$oldQ3.Range("A2").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: add the 2022-Q4 row on top, shift the rest
#    down (2022-Q3 / 2022-Q2 / 2022-Q1 keep their figures, just move a
#    row lower).
# ---------------------------------------------------------------------
# New A5 cell needs the same "s=2" style the other index cells carry -
# grab it from A4 before the values get rewritten.
$totals.Range("A4").Copy()
$totals.Range("A5").PasteSpecial(-4122)

$totals.Range("A5").Value = 3
$totals.Range("B5").Value = "2022-Q1"
$totals.Range("C5").Value = 9
$totals.Range("D5").Value = 0.77

$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2022-Q2"
$totals.Range("C4").Value = 2
$totals.Range("D4").Value = 0

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("C3").Value = 1
$totals.Range("D3").Value = 0.01

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 2
$totals.Range("D2").Value = 0.01

# ---------------------------------------------------------------------
# 3. Restore the originally-active tab: "2022-Q1" was the selected sheet
#    before this edit (and keeps all its data/position as the last
#    sheet), so re-select it now that "2022-Q4" (freshly added/active
#    by default) is done being populated.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
